$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Testing a change on my laptop".
$srcIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Testing a change on my laptop*") {
        $srcIndex = $i
        break
    }
}

if ($srcIndex -eq 0) {
    throw "Could not find paragraph containing 'Testing a change on my laptop'"
}

# The destination paragraph is two paragraphs further down: an empty
# paragraph that follows the empty paragraph right after the source one.
$dstIndex = $srcIndex + 2

# Grab the formatted run (paragraph mark excluded) before we touch anything,
# so the new text later inherits identical run-level formatting.
$srcParaRange = $d.Paragraphs.Item($srcIndex).Range
$srcTextRange = $d.Range($srcParaRange.Start, $srcParaRange.End - 1)
$savedFormattedText = $srcTextRange.FormattedText

# Paste the captured formatting into the destination paragraph first (while
# the saved range reference is still fresh/valid), carrying over the
# original wording for now.
$dstParaRange = $d.Paragraphs.Item($dstIndex).Range
$dstParaRange.FormattedText = $savedFormattedText

# Now remove the run holding "Testing a change on my laptop" from the source
# paragraph entirely, leaving it empty (matches the diff: the <w:r>
# disappears completely). Re-fetch the range fresh since the document
# shifted after the paste above.
$srcParaRange2 = $d.Paragraphs.Item($srcIndex).Range
$srcTextRange2 = $d.Range($srcParaRange2.Start, $srcParaRange2.End - 1)
$srcTextRange2.Delete()

# Finally, overwrite the pasted text in the destination paragraph with the
# new wording while preserving the run formatting that was applied. Re-fetch
# the range fresh once more.
$dstParaRange2 = $d.Paragraphs.Item($dstIndex).Range
$dstTextRange2 = $d.Range($dstParaRange2.Start, $dstParaRange2.End - 1)
$dstTextRange2.Text = "This should change"
